$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns keep their string (inlineStr) representation instead of
# being auto-converted to numbers by Excel when the text looks numeric.
$ws.Range("A2:H9").NumberFormat = "@"
$ws.Range("J2:L9").NumberFormat = "@"
$ws.Range("O2:P9").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "2098"
$ws.Range("B2").Value = "5/24/2024"
$ws.Range("C2").Value = "AZURDUY JUANA 2449"
$ws.Range("D2").Value = "13"
$ws.Range("E2").Value = "788826017"
$ws.Range("F2").Value = "NEW"
$ws.Range("G2").Value = "Pendiente"
$ws.Range("H2").Value = "Terminal con rienda"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = "Cambio"
$ws.Range("K2").Value = "Sin equipos"
$ws.Range("L2").Value = "Pasante"
$ws.Range("M2").Value = -58.467279
$ws.Range("N2").Value = -34.551117
$ws.Range("O2").Value = "Saavedra"
$ws.Range("P2").Value = "Capital Norte"

# Row 3
$ws.Range("A3").Value = "3299"
$ws.Range("B3").Value = "9/10/2024"
$ws.Range("C3").Value = "DIAZ COLODRERO 3309"
$ws.Range("D3").Value = "12"
$ws.Range("E3").Value = "796186684"
$ws.Range("F3").Value = "NEW"
$ws.Range("G3").Value = "Pendiente"
$ws.Range("H3").Value = "Colocar columna para solicitar traspasos"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = "Cambio"
$ws.Range("K3").Value = "Nodo TLC"
$ws.Range("L3").Value = "Pasante"
$ws.Range("M3").Value = -58.491722
$ws.Range("N3").Value = -34.565845
$ws.Range("O3").Value = "Paternal"
$ws.Range("P3").Value = "Capital Norte"

# Row 4
$ws.Range("A4").Value = "3839"
$ws.Range("B4").Value = "10/23/2024"
$ws.Range("C4").Value = "PICO 1511"
$ws.Range("D4").Value = "13"
$ws.Range("E4").Value = "798390296"
$ws.Range("F4").Value = "NEW"
$ws.Range("G4").Value = "Pendiente"
$ws.Range("H4").Value = "Poste inclinado"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = "Aplomo"
$ws.Range("K4").Value = "Sin equipos"
$ws.Range("L4").Value = "Poste"
$ws.Range("M4").Value = -58.465596
$ws.Range("N4").Value = -34.53627
$ws.Range("O4").Value = "Saavedra"
$ws.Range("P4").Value = "Capital Norte"

# Row 5
$ws.Range("A5").Value = "801645368"
$ws.Range("B5").Value = "12/13/2024"
$ws.Range("C5").Value = "San Blas 1809"
$ws.Range("D5").Value = "11"
$ws.Range("E5").Value = "801645368"
$ws.Range("F5").Value = "NEW"
$ws.Range("G5").Value = "Pendiente"
$ws.Range("H5").Value = "Picada"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "Cambio"
$ws.Range("K5").Value = "Sin equipos"
$ws.Range("L5").Value = "Pasante"
$ws.Range("M5").Value = -58.467767
$ws.Range("N5").Value = -34.604588
$ws.Range("O5").Value = "Paternal"
$ws.Range("P5").Value = "Capital Norte"

# Row 6
$ws.Range("A6").Value = "5589"
$ws.Range("B6").Value = "12/31/2023"
$ws.Range("C6").Value = "ARCOS 1520"
$ws.Range("D6").Value = "13"
$ws.Range("E6").Value = "799540526"
$ws.Range("F6").Value = "NEW"
$ws.Range("G6").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H6").Value = "Picada"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = "Pasante"
$ws.Range("M6").Value = -58.449125
$ws.Range("N6").Value = -34.565958
$ws.Range("O6").Value = "Colegiales"
$ws.Range("P6").Value = "Capital Norte"

# Row 7
$ws.Range("A7").Value = "4595"
$ws.Range("B7").Value = "1/15/2025"
$ws.Range("C7").Value = "PAROISSIEN 1806"
$ws.Range("D7").Value = "13"
$ws.Range("E7").Value = "802747617"
$ws.Range("F7").Value = "NEW"
$ws.Range("G7").Value = "Pendiente"
$ws.Range("H7").Value = "Aplomar"
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = "Aplomo"
$ws.Range("K7").Value = "Sin equipos"
$ws.Range("L7").Value = "Terminal"
$ws.Range("M7").Value = -58.464172
$ws.Range("N7").Value = -34.543845
$ws.Range("O7").Value = "Saavedra"
$ws.Range("P7").Value = "Capital Norte"

# Row 8
$ws.Range("A8").Value = "4662"
$ws.Range("B8").Value = "1/21/2025"
$ws.Range("C8").Value = "ALTOLAGUIRRE 2397"
$ws.Range("D8").Value = "12"
$ws.Range("E8").Value = "802823938"
$ws.Range("F8").Value = "NEW"
$ws.Range("G8").Value = "Pendiente"
$ws.Range("H8").Value = "Inclinada"
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = "Aplomo"
$ws.Range("K8").Value = "Sin equipos"
$ws.Range("L8").Value = "Pasante"
$ws.Range("M8").Value = -58.490766
$ws.Range("N8").Value = -34.576987
$ws.Range("O8").Value = "Paternal"
$ws.Range("P8").Value = "Capital Norte"

# Row 9
$ws.Range("A9").Value = "4862"
$ws.Range("B9").Value = "1/23/2025"
$ws.Range("C9").Value = "ARCOS 2263"
$ws.Range("D9").Value = "13"
$ws.Range("E9").Value = "802857379"
$ws.Range("F9").Value = "NEW"
$ws.Range("G9").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H9").Value = "picada"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "Cambio"
$ws.Range("K9").Value = "Nodo Teco"
$ws.Range("L9").Value = "Pasante"
$ws.Range("M9").Value = -58.455082
$ws.Range("N9").Value = -34.558883
$ws.Range("O9").Value = "Saavedra"
$ws.Range("P9").Value = "Capital Norte"
